# Update the "want to go" count (column F) for several events across sheets.
# This mirrors a refreshed data scrape where some counters increased slightly.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1261
$ws1.Range("F5").Value = 5567
$ws1.Range("F6").Value = 1785
$ws1.Range("F9").Value = 1919
$ws1.Range("F16").Value = 49
$ws1.Range("F17").Value = 7909
$ws1.Range("F18").Value = 7909
$ws1.Range("F25").Value = 2
$ws1.Range("F30").Value = 1748

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9540

# Sheet "全部类型" (All Types) - aggregated view of all rows above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9540
$ws4.Range("F5").Value = 1261
$ws4.Range("F10").Value = 5567
$ws4.Range("F12").Value = 1785
$ws4.Range("F15").Value = 1919
$ws4.Range("F22").Value = 49
$ws4.Range("F23").Value = 7909
$ws4.Range("F24").Value = 7909
$ws4.Range("F31").Value = 2
$ws4.Range("F35").Value = 1748
